$d = $word.ActiveDocument

# The "Ancillary Structures" page used to show the substation setback
# diagram as an embedded picture. It is replaced by a plain hyperlink
# that simply displays (and points at) the image's URL on ura.gov.sg.
$imageUrl = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/B102_Setbacks_Ancillary_Structures_Substation.jpg?h=100%25&w=100%25"

if ($d.InlineShapes.Count -ge 1) {
    $shp = $d.InlineShapes.Item(1)
    $r = $shp.Range
    $shp.Delete()
    $d.Hyperlinks.Add($r, $imageUrl, [Type]::Missing, [Type]::Missing, $imageUrl, [Type]::Missing)
}
